# SAM TODO.xlsx - TODO list update from 10/23/14 meeting
$wb = $excel.ActiveWorkbook
$todo = $wb.Worksheets.Item("To Do")
$ideas = $wb.Worksheets.Item("Project Ideas")

# ---- "Project Ideas" sheet: new idea added at the bottom ----
$ideas.Range("A11").Value = "ITC over multiple years"

# ---- "To Do" sheet: status column (A) updates ----
# Several items that were "Not done" are now "Done"
$todo.Range("A36").Value = "Done"
$todo.Range("A38").Value = "Done"
$todo.Range("A39").Value = "Done"
$todo.Range("A58").Value = "Done"
$todo.Range("A63").Value = "Done"
$todo.Range("A64").Value = "Done"
$todo.Range("A66").Value = "Done"
$todo.Range("A111").Value = "Done"
$todo.Range("A112").Value = "Done"

# Several items that were "Not done" are now "Future"
$todo.Range("A49").Value = "Future"
$todo.Range("A65").Value = "Future"
$todo.Range("A67").Value = "Future"
$todo.Range("A68").Value = "Future"
$todo.Range("A74").Value = "Future"

# Two items get a new, one-off status of "POUT"
$todo.Range("A45").Value = "POUT"
$todo.Range("A46").Value = "POUT"

# Fill in missing "Who" (column C) for a couple of rows
$todo.Range("C52").Value = "Steve"
$todo.Range("C60").Value = "Janine"

# Newly added status markers for rows that previously had no status set
$todo.Range("A123").Value = "Done"
$todo.Range("A124").Value = "Not done"
$todo.Range("A125").Value = "Not done"

# ---- View state: "To Do" tab becomes the active/selected sheet ----
$ideas.Activate()
$ideas.Range("A12").Select()

$todo.Activate()
$todo.Range("A125").Select()
